$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates are Excel serial numbers, same style as existing date column)
$newRows = @(
    @{ Row = 230; A = 44304; B = 3; C = 5; D = 94.6969696969697 },
    @{ Row = 231; A = 44305; B = 2; C = 7; D = 132.5757575757576 },
    @{ Row = 232; A = 44306; B = 0; C = 7; D = 132.5757575757576 },
    @{ Row = 233; A = 44307; B = 0; C = 7; D = 132.5757575757576 }
)

# Copy the formatting of the last existing data row (229) so the new
# date cells in column A keep the same style (centered, bordered, date number format)
$ws.Range("A229").Copy() | Out-Null

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("A$row").Value = $r.A

    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
}

$excel.CutCopyMode = 0
